# Roster order changed: D'Angelo Russell moved up to row 13 (No. 1),
# swapping places with Mo Bamba who moves down to row 14 (No. 12).
# Column A ("No." - the 0-based row index) stays as-is; every other
# column (B..K: uniform No., Player, Pos, Ht, Wt, Birth Date,
# nationality flag, Exp, College, bbref url) is swapped between the
# two rows.
#
# We swap via Copy (instead of reading/writing .Value, whose getter is
# unusable here and whose setter auto-converts numeric-looking text
# such as the Exp column's "4"/"7" into real numbers) so that cell
# types/formatting are preserved exactly, using a scratch range as
# temporary holding space, then clearing that scratch range afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 13
$row2 = 14
$scratchRow = 100

$rangeRow1 = $ws.Range("B$row1`:K$row1")
$rangeRow2 = $ws.Range("B$row2`:K$row2")
$scratch = $ws.Range("B$scratchRow`:K$scratchRow")

$rangeRow1.Copy($scratch)
$rangeRow2.Copy($rangeRow1)
$scratch.Copy($rangeRow2)
$scratch.Clear()
